$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Beneficiarios")

$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0
$ws.Range("I2").Value = "SI"
$ws.Range("K2").Value = "NO"

$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0
$ws.Range("I3").Value = "SI"
$ws.Range("K3").Value = "NO"

$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0
$ws.Range("I4").Value = "SI"
$ws.Range("K4").Value = "NO"

$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0
$ws.Range("I5").Value = "SI"

$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 0
$ws.Range("I6").Value = "SI"

$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0
$ws.Range("I7").Value = "SI"

$ws.Range("F8").Value = 0
$ws.Range("I8").Value = "SI"

$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0
$ws.Range("I9").Value = "SI"

$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 0
$ws.Range("I10").Value = "SI"

$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 0
$ws.Range("I11").Value = "SI"

$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 0
$ws.Range("I12").Value = "SI"

$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 0
$ws.Range("I13").Value = "SI"

$ws.Range("F14").Value = 0
$ws.Range("I14").Value = "SI"

$ws.Range("F15").Value = 0
$ws.Range("I15").Value = "SI"

$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 0
$ws.Range("I16").Value = "SI"

$ws.Range("F17").Value = 0
$ws.Range("I17").Value = "SI"

$ws.Range("F18").Value = 0
$ws.Range("I18").Value = "SI"

$ws.Range("F19").Value = 0
$ws.Range("I19").Value = "SI"

$ws.Range("F20").Value = 0
$ws.Range("I20").Value = "SI"

$ws.Range("F21").Value = 0
$ws.Range("I21").Value = "SI"

$ws.Range("E22").Value = 2
$ws.Range("F22").Value = 0
$ws.Range("I22").Value = "SI"

$ws.Range("E23").Value = 2
$ws.Range("F23").Value = 0
$ws.Range("I23").Value = "SI"

$ws.Range("E24").Value = 2
$ws.Range("F24").Value = 0
$ws.Range("I24").Value = "SI"

$ws.Range("F25").Value = 0
$ws.Range("I25").Value = "SI"

$ws.Range("F26").Value = 0
$ws.Range("I26").Value = "SI"
